# The post that used to occupy row 580 ("「夢みれば叶う」...") was removed
# from the source data. Delete that entire row; Excel will shift all
# subsequent rows up by one and the used range naturally becomes A1:C690.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(580).Delete()
